$wb = $excel.ActiveWorkbook

# --- Update the conversion text on sheet "Hoja1" (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.38 = 49727.55 pesos`n✅ 49727.55 pesos = 12.3 = 976.84 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 80.75
$ws2.Range("O10").Value = 4015.5

$ws2.Range("N12").Value = 4042
$ws2.Range("O12").Value = 79.40000000000001
